# Update meter number values for the two bulk-imported space rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "M-103"
$ws.Range("J3").Value = "M-104"

# Leave the selection where the user's edit session ended up
$ws.Range("B4").Select()
